$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.391.79'
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").Value = '3.359.33'
$ws.Range("E3").Value = '  -1.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.79%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.358.30'
$ws.Range("E8").Value = '  -1.64%  '

$ws.Range("E9").Value = '  -1.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.97%  '

$ws.Range("E12").Value = '  -1.43%  '

$ws.Range("D13").Value = '3.931.31'
$ws.Range("E13").Value = '  -1.65%  '

$ws.Range("E14").Value = '  +1.63%  '

$ws.Range("E15").Value = '  -4.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.28%  '

$ws.Range("D17").Value = '3.361.96'
$ws.Range("E17").Value = '  -1.67%  '

$ws.Range("D18").Value = '61.438.52'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.16%  '

$ws.Range("E20").Value = '  -1.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.32'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '375.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.98%  '

$ws.Range("E23").Value = '  -3.89%  '

$ws.Range("D24").Value = '3.501.44'
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.01%  '

$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.15%  '

$ws.Range("E30").Value = '  +0.45%  '

$ws.Range("E31").Value = '  +2.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.93%  '

$ws.Range("E33").Value = '  -0.79%  '

$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.35%  '

$ws.Range("E36").Value = '  -6.12%  '

$ws.Range("E37").Value = '  -3.69%  '

$ws.Range("E38").Value = '  -2.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0765'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.43%  '

$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("E42").Value = '  -1.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.771'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.69%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.14%  '

$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.84'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").Value = '2.355.42'
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("E51").Value = '  -2.22%  '
